# Insert a new weekly record for "Feria Lagunitas de Puerto Montt - Lechuga".
# This shifts rows 856:882 down to 857:883 and populates the new row 856
# with the latest observation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 856 (existing rows 856-882 shift to 857-883).
$ws.Rows.Item(856).Insert()

# Populate the newly inserted row 856 with the new weekly record.
$ws.Cells.Item(856, 1).Value = 4
$ws.Cells.Item(856, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(856, 3).Value = "Los Lagos"
$ws.Cells.Item(856, 4).Value = "2023-05-29"
$ws.Cells.Item(856, 5).Value = 10
$ws.Cells.Item(856, 6).Value = 100112033
$ws.Cells.Item(856, 7).Value = "Lechuga"
$ws.Cells.Item(856, 8).Value = "Escarola"
$ws.Cells.Item(856, 9).Value = "Primera"
$ws.Cells.Item(856, 10).Value = 150
$ws.Cells.Item(856, 11).Value = 10500
$ws.Cells.Item(856, 12).Value = 10500
$ws.Cells.Item(856, 13).Value = 10500
$ws.Cells.Item(856, 14).Value = "`$/caja 15 unidades"
$ws.Cells.Item(856, 15).Value = "Región de Coquimbo"
$ws.Cells.Item(856, 16).Value = 700
$ws.Cells.Item(856, 17).Value = 15
$ws.Cells.Item(856, 18).Value = "Hortaliza"
